# Update calculated values in pl_mw.xlsx Sheet1 for the 380 kV case
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.466837306709067
    "C2" = 0.0398775350072782
    "D2" = 0.1873939752032499
    "E2" = 0.166210065083483
    "F2" = 1.484535831506925
    "J2" = 0.1890045645124232
    "K2" = 0.429863010182487
    "N2" = 1.667300998228932
    "O2" = 3.623486681340438
    "B3" = 0.4287475654596165
    "C3" = 0.03478104256780057
    "D3" = 0.1815751677574298
    "E3" = 0.1620237846593753
    "F3" = 1.485436559375884
    "J3" = 0.184997984565868
    "K3" = 0.3893307655359308
    "N3" = 1.684937734256404
    "O3" = 3.638125709781917
    "B4" = 0.4054897257903463
    "C4" = 0.03164125845404442
    "D4" = 0.1780815867062699
    "E4" = 0.1595368129236796
    "F4" = 1.486746116947181
    "J4" = 0.1826434735644256
    "K4" = 0.3645219074932697
    "N4" = 1.696319255513841
    "O4" = 3.648967928411508
    "B5" = 0.3960450350762699
    "C5" = 0.03035919039118085
    "D5" = 0.1766779432699508
    "E5" = 0.1585443799780357
    "F5" = 1.487470109109765
    "J5" = 0.1817105536292871
    "K5" = 0.3544322752698292
    "N5" = 1.701096241540294
    "O5" = 3.653852443779911
    "B6" = 0.3944787640630807
    "C6" = 0.03014614999635512
    "D6" = 0.1764460815811617
    "E6" = 0.1583808588063071
    "F6" = 1.487601826499009
    "J6" = 0.1815572479465573
    "K6" = 0.3527581341589325
    "N6" = 1.701897846966586
    "O6" = 3.654691674725214
    "B7" = 0.4053622168017057
    "C7" = 0.03162397839965081
    "D7" = 0.1780625754816469
    "E7" = 0.15952334341506
    "F7" = 1.486755110151257
    "J7" = 0.1826307842997466
    "K7" = 0.364385752766367
    "N7" = 1.696383117128624
    "O7" = 3.649031914882983
    "B8" = 0.4536774065601321
    "C8" = 0.03812248416728892
    "D8" = 0.1853712760129866
    "E8" = 0.1647493505914639
    "F8" = 1.484689449904295
    "J8" = 0.1876012014553581
    "K8" = 0.4158715922706051
    "N8" = 1.67326747299223
    "O8" = 3.628149500113125
    "B9" = 0.549432037787227
    "C9" = 0.05078036546194653
    "D9" = 0.200328170828584
    "E9" = 0.1756578958780466
    "F9" = 1.486638268594092
    "J9" = 0.1981853946697782
    "K9" = 0.5174363035352769
    "N9" = 1.632321110503864
    "O9" = 3.601907847604792
    "B10" = 0.6203807254739786
    "C10" = 0.06002577945206156
    "D10" = 0.2116939613100328
    "E10" = 0.184073931529646
    "F10" = 1.491726132786326
    "J10" = 0.2064727756210232
    "K10" = 0.5924049034930761
    "N10" = 1.604908412260271
    "O10" = 3.591598606151763
    "B11" = 0.652783708896294
    "C11" = 0.06421956658071792
    "D11" = 0.2169456804100776
    "E11" = 0.1879896311423366
    "F11" = 1.494834492637878
    "J11" = 0.2103541482696585
    "K11" = 0.6265825886491712
    "N11" = 1.593017225745982
    "O11" = 3.58885764190461
    "B12" = 0.665071854851476
    "C12" = 0.06580586585292281
    "D12" = 0.218945979398228
    "E12" = 0.1894849071571585
    "F12" = 1.496125647310592
    "J12" = 0.2118399368558528
    "K12" = 0.6395349956356142
    "N12" = 1.588597615513275
    "O12" = 3.588099970784299
    "B13" = 0.6624245985353525
    "C13" = 0.06546430897910227
    "D13" = 0.2185146656916572
    "E13" = 0.1891623182800259
    "F13" = 1.495842501159544
    "J13" = 0.2115192345895167
    "K13" = 0.6367450251275955
    "N13" = 1.589545750708316
    "O13" = 3.588250682931459
    "B14" = 0.6537943078373019
    "C14" = 0.06435010879853564
    "D14" = 0.2171100146406246
    "E14" = 0.188112398481934
    "F14" = 1.494938430656617
    "J14" = 0.2104760645167119
    "K14" = 0.6276479927722391
    "N14" = 1.592651951424741
    "O14" = 3.588789690519235
    "B15" = 0.6485103116750111
    "C15" = 0.06366739288553447
    "D15" = 0.2162511309136477
    "E15" = 0.1874709166665269
    "F15" = 1.494399516668281
    "J15" = 0.2098391753367537
    "K15" = 0.6220770901435912
    "N15" = 1.59456544231844
    "O15" = 3.589156349269643
    "B16" = 0.6182656420854755
    "C16" = 0.05975145804798387
    "D16" = 0.2113523767058183
    "E16" = 0.18381978147643
    "F16" = 1.491538962727248
    "J16" = 0.2062213585673476
    "K16" = 0.5901727514633706
    "N16" = 1.605697185576932
    "O16" = 3.591816952347472
    "B17" = 0.5997439082712788
    "C17" = 0.05734603491080748
    "D17" = 0.2083679085380936
    "E17" = 0.1816022246264168
    "F17" = 1.489987383283307
    "J17" = 0.2040304619397801
    "K17" = 0.5706190266997453
    "N17" = 1.612674488373436
    "O17" = 3.593948301275049
    "B18" = 0.5891027858824032
    "C18" = 0.05596137410451263
    "D18" = 0.2066589865260511
    "E18" = 0.1803349557725227
    "F18" = 1.489169682377977
    "J18" = 0.2027808031476326
    "K18" = 0.5593792509651223
    "N18" = 1.616742141277443
    "O18" = 3.595357620847864
    "B19" = 0.5855019785529407
    "C19" = 0.05549236070316965
    "D19" = 0.2060816948336992
    "E19" = 0.1799072922367344
    "F19" = 1.488905659406797
    "J19" = 0.20235949198063
    "K19" = 0.5555748780585645
    "N19" = 1.61812873589201
    "O19" = 3.595866294196554
    "B20" = 0.6017143314647058
    "C20" = 0.05760221326735859
    "D20" = 0.2086848177748379
    "E20" = 0.1818374379024164
    "F20" = 1.49014481815469
    "J20" = 0.2042626014225846
    "K20" = 0.5726998329275546
    "N20" = 1.611926102297032
    "O20" = 3.593702431526793
    "B21" = 0.6563287538614873
    "C21" = 0.06467742573747159
    "D21" = 0.2175222811386703
    "E21" = 0.1884204469554334
    "F21" = 1.495200882373851
    "J21" = 0.2107820348803102
    "K21" = 0.6303197434259857
    "N21" = 1.591737323045725
    "O21" = 3.588623764064579
    "B22" = 0.6921261495429007
    "C22" = 0.06929097923779182
    "D22" = 0.2233655592518744
    "E22" = 0.192795573955415
    "F22" = 1.499170261686785
    "J22" = 0.215136093617005
    "K22" = 0.6680360431938084
    "N22" = 1.579028497735139
    "O22" = 3.586938224310842
    "B23" = 0.6730111066834183
    "C23" = 0.06682962392098091
    "D23" = 0.2202407543040437
    "E23" = 0.1904538476578281
    "F23" = 1.49699090974957
    "J23" = 0.2128037272823917
    "K23" = 0.6479010011782691
    "N23" = 1.585766971922476
    "O23" = 3.587688334067735
    "B24" = 0.6008234804529877
    "C24" = 0.05748640048206255
    "D24" = 0.2085415216443636
    "E24" = 0.1817310741976925
    "F24" = 1.490073410332343
    "J24" = 0.2041576202716868
    "K24" = 0.5717590944033759
    "N24" = 1.612264272245757
    "O24" = 3.593813016155366
    "B25" = 0.5234216593057113
    "C25" = 0.04736546256530971
    "D25" = 0.1962154622055152
    "E25" = 0.1726362563147177
    "F25" = 1.485468914126898
    "J25" = 0.1952323776751967
    "K25" = 0.4898978073256046
    "N25" = 1.642929097180112
    "O25" = 3.607431865002297
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

